$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1113
$ws1.Range("G3").Value = "不可售"
$ws1.Range("F4").Value = 1525
$ws1.Range("F5").Value = 8814
$ws1.Range("F9").Value = 304
$ws1.Range("F10").Value = 162
$ws1.Range("F11").Value = 24
$ws1.Range("F12").Value = 22
$ws1.Range("F13").Value = 3674
$ws1.Range("F17").Value = 2920
$ws1.Range("F18").Value = 152
$ws1.Range("F21").Value = 217
$ws1.Range("F22").Value = 2463
$ws1.Range("F23").Value = 82

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1113
$ws4.Range("G3").Value = "不可售"
$ws4.Range("F4").Value = 1525
$ws4.Range("F5").Value = 8814
$ws4.Range("F9").Value = 304
$ws4.Range("F10").Value = 162
$ws4.Range("F11").Value = 24
$ws4.Range("F12").Value = 22
$ws4.Range("F13").Value = 3674
$ws4.Range("F17").Value = 2920
$ws4.Range("F18").Value = 152
$ws4.Range("F21").Value = 217
$ws4.Range("F22").Value = 2463
$ws4.Range("F24").Value = 82

$wb.Save()
